$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Fill in the "Method Inputs" / "Expected Result" test-plan data.
# Cells are written in the same order the author originally typed them so
# that the shared-string table indices line up with the target workbook.
# ---------------------------------------------------------------------------

# Row 7 - __init__ / Attributes set to input values.
$ws.Range("E7").Value = "None"
$ws.Range("G7").Value = "Attribute got set"
$ws.Range("F7").Value = "client id:525" + "`n" + "first name :Daniel" + "`n" + "last name:Chinchilla" + "`n" + "email address: erichilla0525@gmail.com"

# Row 12 - client_number / Returns client_number attribute.
$ws.Range("F12").Value = "client id:525"

# Row 13 - first_name / Returns first_name attribute.
$ws.Range("F13").Value = "first name :Daniel"

# Row 14 - last_name / Returns last_name attribute.
$ws.Range("F14").Value = "last name:Chinchilla"

# Row 15 - email_address / Returns email_address attribute.
$ws.Range("F15").Value = "email address: erichilla0525@gmail.com"

$ws.Range("G13").Value = "Daniel"
$ws.Range("G14").Value = "Chinchilla"

# G15 becomes a mailto: hyperlink displaying the raw address.
$ws.Range("G15").Value = "erichilla0525@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G15"), "mailto:erichilla0525@gmail.com")

# Row 8 - Exception raised when invalid client number.
$ws.Range("G8").Value = "valueerror rasies"
$ws.Range("F8").Value = "client id: invalid client number" + "`n" + "first name :Daniel" + "`n" + "last name:Chinchilla" + "`n" + "email address: erichilla0525@gmail.com"

# Row 9 - Exception raised when blank first_name.
$ws.Range("F9").Value = "client id:525" + "`n" + "first name :invalid first name" + "`n" + "last name:Chinchilla" + "`n" + "email address: erichilla0525@gmail.com"

# Row 10 - Exception raised when blank last_name.
$ws.Range("F10").Value = "client id:525" + "`n" + "first name :Daniel" + "`n" + "last name:invalid last name" + "`n" + "email address: erichilla0525@gmail.com"

# Row 11 - Email address set to default value when invalid.
$ws.Range("F11").Value = "client id:525" + "`n" + "first name :Daniel" + "`n" + "last name:Chinchilla" + "`n" + "email address: invalid email address"

# Row 16 - __str__ / Returns string in expected format.
$ws.Range("G16").Value = "Chinchilla, Daniel [525] - email@pixell-river.com"

# Remaining cells reuse already-introduced strings.
$ws.Range("E8").Value = "None"
$ws.Range("E9").Value = "None"
$ws.Range("E10").Value = "None"
$ws.Range("E11").Value = "None"
$ws.Range("E12").Value = "None"
$ws.Range("E13").Value = "None"
$ws.Range("E14").Value = "None"
$ws.Range("E15").Value = "None"

$ws.Range("G9").Value = "valueerror rasies"
$ws.Range("G10").Value = "valueerror rasies"
$ws.Range("G11").Value = "valueerror rasies"

$ws.Range("F16").Value = "client id:525" + "`n" + "first name :Daniel" + "`n" + "last name:Chinchilla" + "`n" + "email address: erichilla0525@gmail.com"
$ws.Rows("16:16").RowHeight = 76.5

# G12 is the numeric client id, bold (inherited from the row style) and
# explicitly left aligned.
$ws.Range("G12").Value = 525
$ws.Range("G12").HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# View state: scrolled down a bit with E10 selected.
# ---------------------------------------------------------------------------
$ws.Range("E10").Select()
